$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New set of accuracy values (percentages) once the sample count grows from 270 to 271
# (correct_count / 271 * 100), including a brand-new 50th row.
$values = @(
    71.58671586715867,
    71.58671586715867,
    71.58671586715867,
    71.21771217712177,
    71.58671586715867,
    71.58671586715867,
    73.43173431734317,
    73.80073800738008,
    74.53874538745387,
    73.80073800738008,
    72.32472324723247,
    73.06273062730627,
    74.53874538745387,
    74.53874538745387,
    74.53874538745387,
    74.16974169741697,
    74.16974169741697,
    74.53874538745387,
    73.06273062730627,
    73.43173431734317,
    73.80073800738008,
    70.11070110701108,
    70.11070110701108,
    69.74169741697416,
    73.06273062730627,
    72.69372693726937,
    73.80073800738008,
    74.90774907749078,
    74.53874538745387,
    74.16974169741697,
    69.74169741697416,
    70.11070110701108,
    71.58671586715867,
    70.47970479704797,
    70.47970479704797,
    75.64575645756457,
    69.00369003690037,
    69.37269372693727,
    70.47970479704797,
    74.16974169741697,
    74.16974169741697,
    74.90774907749078,
    74.16974169741697,
    70.47970479704797,
    71.21771217712177,
    70.47970479704797,
    74.53874538745387,
    70.84870848708486,
    73.43173431734317
)

# Row 1 holds the dataset label, repeated as the header for the column.
$ws.Range("A1").Value = "HK_G_acc_G"

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
